$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two day headers and add a third day column
$ws.Range("C2").Value = "Tâches jour 1 "
$ws.Range("D2").Value = "Tâches  jour 2"
$ws.Range("E2").Value = "Tâches  jour 3"

# New column E width (characters) -- closest reproducible value to the
# saved workbook's 34.88671875 given the engine's pixel-grid quantisation
$ws.Columns.Item(5).ColumnWidth = 34

# Copy the border/style formatting from column D into column E for the used rows
$ws.Range("D2:D5").Copy() | Out-Null
$ws.Range("E2:E5").PasteSpecial(-4122) | Out-Null

# Fill in the new "jour 3" column values (order matters for shared-string ids)
$ws.Range("E5").Value = "Back Controller (partie admin)+Front"
$ws.Range("E3").Value = "Back Controller (partie admin)+front"
$ws.Range("E4").Value = "Back Controller (partie responsable)+front"

# Move the active selection to E7, matching the saved workbook state
$ws.Range("E7").Select() | Out-Null
